$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Word's Range.Text for a table cell includes the trailing cell-mark
# characters (CR + BEL); strip those (plus ordinary whitespace) before
# comparing/testing cell contents.
function CellText($cell) {
    return $cell.Range.Text.TrimEnd([char]13, [char]7).Trim()
}

# Locate the "K = 2" results column from the header row so we don't rely on
# a hard-coded column index.
$headerRow = 1
$k2Col = 0
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    if ((CellText $t.Cell($headerRow, $c)) -eq "K = 2") {
        $k2Col = $c
    }
}

# New k=2 spatial conStruct admixture proportions for this table's rows,
# in top-to-bottom order, filling only the still-blank cells.
$values = @("0.05", "0.95", "0.00", "0.00", "0.00")

$valueIndex = 0
for ($r = $headerRow + 1; $r -le $t.Rows.Count; $r++) {
    $cell = $t.Cell($r, $k2Col)
    if ((CellText $cell) -eq "") {
        $cell.Range.Text = $values[$valueIndex]
        $valueIndex++
    }
}
